# Auto-generated Excel COM-interop edit script
# Applies the diff between before.xlsx and the target (after) workbook:
#  - "展览" (Exhibitions) sheet: refresh "想去人数" (F) counts; one price (G3) now sold out
#  - "演出" (Performances) sheet: refresh one "想去人数" (F5) count
#  - "全部类型" (All types) sheet: the aggregated feed was regenerated - one listing
#    ("胡良伟专场") dropped out and rows 3-18 shifted up by one to absorb the new
#    listing ("钟晨瑶内场票") that now appears at the end of the block; every F number
#    downstream was refreshed to the latest scrape.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 2711
$ws.Cells.Item(3, 7).Value = "已售罄"
$ws.Cells.Item(5, 6).Value = 1519
$ws.Cells.Item(6, 6).Value = 1140
$ws.Cells.Item(9, 6).Value = 1162
$ws.Cells.Item(11, 6).Value = 118
$ws.Cells.Item(12, 6).Value = 552
$ws.Cells.Item(13, 6).Value = 9202
$ws.Cells.Item(14, 6).Value = 399
$ws.Cells.Item(15, 6).Value = 2501
$ws.Cells.Item(16, 6).Value = 7
$ws.Cells.Item(17, 6).Value = 259
$ws.Cells.Item(18, 6).Value = 182
$ws.Cells.Item(20, 6).Value = 632
$ws.Cells.Item(23, 6).Value = 999
$ws.Cells.Item(24, 6).Value = 2095
$ws.Cells.Item(25, 6).Value = 2189
$ws.Cells.Item(27, 6).Value = 1890
$ws.Cells.Item(28, 6).Value = 1929
$ws.Cells.Item(29, 6).Value = 482
$ws.Cells.Item(31, 6).Value = 277
$ws.Cells.Item(32, 6).Value = 159
$ws.Cells.Item(33, 6).Value = 212
$ws.Cells.Item(35, 6).Value = 325
$ws.Cells.Item(36, 6).Value = 63
$ws.Cells.Item(37, 6).Value = 296
$ws.Cells.Item(38, 6).Value = 490
$ws.Cells.Item(39, 6).Value = 11
$ws.Cells.Item(40, 6).Value = 60
$ws.Cells.Item(41, 6).Value = 578
$ws.Cells.Item(42, 6).Value = 35
$ws.Cells.Item(43, 6).Value = 1396
$ws.Cells.Item(44, 6).Value = 302
$ws.Cells.Item(45, 6).Value = 9
$ws.Cells.Item(46, 6).Value = 176
$ws.Cells.Item(47, 6).Value = 642
$ws.Cells.Item(49, 6).Value = 297

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 24

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 2711
$ws.Cells.Item(3, 3).Value = "杭州·OZ·富坚义博only"
$ws.Cells.Item(3, 4).Value = "北干街道萧杭路689号浙农东巢艺术公园 Fashion Bund时尚外滩艺术中心"
$ws.Cells.Item(3, 5).Value = "2024.03.16 10:00-03.16 17:00"
$ws.Cells.Item(3, 6).Value = 355
$ws.Cells.Item(3, 7).Value = 88
$ws.Cells.Item(3, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81151"
$ws.Cells.Item(3, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/CxqdeAPa1705658329588.jpeg"
$ws.Cells.Item(4, 3).Value = "杭州·SST动漫嘉年华"
$ws.Cells.Item(4, 4).Value = "沈半路171号 Tcar汽车文化主题公园"
$ws.Cells.Item(4, 5).Value = "2024.03.16 09:00-03.17 17:00"
$ws.Cells.Item(4, 6).Value = 1519
$ws.Cells.Item(4, 7).Value = 68
$ws.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81196"
$ws.Cells.Item(4, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/yFyT2uAT1705897787652.jpeg"
$ws.Cells.Item(5, 3).Value = "杭州·《挪威的森林》摇滚情歌之夜--630乐团演绎经典"
$ws.Cells.Item(5, 4).Value = "湖墅南路136-138号 浙话艺术剧院"
$ws.Cells.Item(5, 5).Value = "2024.03.16 19:00-03.16 21:00"
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 112
$ws.Cells.Item(5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81557"
$ws.Cells.Item(5, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/NXR7ATah1706682091721.jpeg"
$ws.Cells.Item(6, 3).Value = "杭州·排球少年*蓝锁ONLY"
$ws.Cells.Item(6, 4).Value = "亚太路湘湖3期东南侧约290米 原创壹号羽毛球馆"
$ws.Cells.Item(6, 5).Value = "2024.03.16 10:00-03.16 17:00"
$ws.Cells.Item(6, 6).Value = 1140
$ws.Cells.Item(6, 7).Value = 60
$ws.Cells.Item(6, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81075"
$ws.Cells.Item(6, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/9AL6kYuj1705634962275.jpeg"
$ws.Cells.Item(7, 3).Value = "杭州·春和景明代号鸢only"
$ws.Cells.Item(7, 4).Value = "金沙大道681号 金沙湖大剧院"
$ws.Cells.Item(7, 5).Value = "2024.03.16 09:30-03.16 16:00"
$ws.Cells.Item(7, 7).Value = "已售罄"
$ws.Cells.Item(7, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81894"
$ws.Cells.Item(7, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/1RX6wnFN1708323470900.png"
$ws.Cells.Item(8, 3).Value = "杭州·百鬼夜行·咒术回战only"
$ws.Cells.Item(8, 4).Value = "长生路18号 梅地亚宾馆"
$ws.Cells.Item(8, 5).Value = "2024.03.16 09:00-03.16 17:00"
$ws.Cells.Item(8, 7).Value = 79
$ws.Cells.Item(8, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81478"
$ws.Cells.Item(8, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/4weHdCdk1706495040356.jpeg"
$ws.Cells.Item(9, 3).Value = "杭州·造梦探险家——次元茶话会"
$ws.Cells.Item(9, 4).Value = "临平街道北沙西路156-1号 杭州临平遇上设计师酒店"
$ws.Cells.Item(9, 5).Value = "2024.03.16 10:00-03.16 17:00"
$ws.Cells.Item(9, 6).Value = 118
$ws.Cells.Item(9, 7).Value = 38
$ws.Cells.Item(9, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81914"
$ws.Cells.Item(9, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/HHHVcvcC1709015213282.png"
$ws.Cells.Item(10, 2).NumberFormat = "@"
$ws.Cells.Item(10, 2).Value = "2024-03-23"
$ws.Cells.Item(10, 3).Value = "杭州·AD02动漫展"
$ws.Cells.Item(10, 4).Value = "浙江省杭州市萧山区奔竞大道353号 国际博览中心"
$ws.Cells.Item(10, 5).Value = "2024.03.23 10:00-03.24 17:00"
$ws.Cells.Item(10, 6).Value = 9202
$ws.Cells.Item(10, 7).Value = 75
$ws.Cells.Item(10, 8).Value = "https://show.bilibili.com/platform/detail.html?id=80905"
$ws.Cells.Item(10, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/D3QaPamg1705397424553.jpeg"
$ws.Cells.Item(11, 3).Value = "杭州·AD02动漫展  青柳尊哉内场票"
$ws.Cells.Item(11, 5).Value = "2024.03.23 10:00-03.23 17:00"
$ws.Cells.Item(11, 6).Value = 399
$ws.Cells.Item(11, 7).Value = 528
$ws.Cells.Item(11, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81503"
$ws.Cells.Item(11, 9).Value = "//i1.hdslb.com/bfs/openplatform/202401/OmqxboDC1706522627528.jpeg"
$ws.Cells.Item(12, 3).Value = "杭州·AD02动漫展--卡琳娜签售票"
$ws.Cells.Item(12, 5).Value = "2024.03.23 09:30-03.23 17:00"
$ws.Cells.Item(12, 7).Value = "已售罄"
$ws.Cells.Item(12, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81941"
$ws.Cells.Item(12, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/CZjxY9ZC1708416661613.jpeg"
$ws.Cells.Item(13, 3).Value = "杭州·《天空之城》久石让·宫崎骏动漫经典作品音乐会|浙江电影爱乐乐团"
$ws.Cells.Item(13, 4).Value = "武林路77号 文化馆小剧场"
$ws.Cells.Item(13, 5).Value = "2024.03.23 19:30-03.23 21:00"
$ws.Cells.Item(13, 7).Value = 90
$ws.Cells.Item(13, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82065"
$ws.Cells.Item(13, 9).Value = "//i0.hdslb.com/bfs/openplatform/202402/K7MwIOqE1708918668985.jpeg"
$ws.Cells.Item(14, 3).Value = "杭州·星玫Rostar偶像团 1st off会 - 莫里生日SP"
$ws.Cells.Item(14, 4).Value = "下沙大道30号 杭州璞砚酒店"
$ws.Cells.Item(14, 5).Value = "2024.03.23 12:00-03.23 21:00"
$ws.Cells.Item(14, 6).Value = 7
$ws.Cells.Item(14, 7).Value = 58
$ws.Cells.Item(14, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82690"
$ws.Cells.Item(14, 9).Value = "//i2.hdslb.com/bfs/openplatform/202403/grcA9fYK1710327154137.jpeg"
$ws.Cells.Item(15, 3).Value = "浙江·燃爆全场·世界电影主题音乐会 《复仇者联盟》、《歌剧魅影》、《泰坦尼克号》燃情主题音乐"
$ws.Cells.Item(15, 4).Value = "曙光路31号 浙江音乐厅"
$ws.Cells.Item(15, 5).Value = "2024.03.23 19:30-03.23 21:00"
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 100
$ws.Cells.Item(15, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82685"
$ws.Cells.Item(15, 9).Value = "//i0.hdslb.com/bfs/openplatform/202403/KFRQDTnB1710210073027.jpeg"
$ws.Cells.Item(16, 2).NumberFormat = "@"
$ws.Cells.Item(16, 2).Value = "2024-03-24"
$ws.Cells.Item(16, 3).Value = "杭州·AD02动漫展  岩永彻也内场票"
$ws.Cells.Item(16, 4).Value = "浙江省杭州市萧山区奔竞大道353号 国际博览中心"
$ws.Cells.Item(16, 5).Value = "2024.03.24 10:00-03.24 17:00"
$ws.Cells.Item(16, 6).Value = 259
$ws.Cells.Item(16, 7).Value = 528
$ws.Cells.Item(16, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81239"
$ws.Cells.Item(16, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/hww9WUpD1705914756383.jpeg"
$ws.Cells.Item(17, 3).Value = "杭州·AD02动漫展--亦之紫F、L句号内场票"
$ws.Cells.Item(17, 4).Value = "钱江世纪城奔竞大道353号 杭州国际博览中心"
$ws.Cells.Item(17, 5).Value = "2024.03.24 12:00-03.24 16:00"
$ws.Cells.Item(17, 6).Value = 182
$ws.Cells.Item(17, 7).Value = 258
$ws.Cells.Item(17, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81836"
$ws.Cells.Item(17, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/ecrRfQce1707375167618.jpeg"
$ws.Cells.Item(18, 3).Value = "杭州·AD02动漫展--钟晨瑶内场票"
$ws.Cells.Item(18, 5).Value = "2024.03.24 09:30-03.24 17:00"
$ws.Cells.Item(18, 7).Value = "已售罄"
$ws.Cells.Item(18, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81820"
$ws.Cells.Item(18, 9).Value = "//i1.hdslb.com/bfs/openplatform/202402/aHRmCxr31707296105225.jpeg"
$ws.Cells.Item(19, 6).Value = 632
$ws.Cells.Item(21, 6).Value = 999
$ws.Cells.Item(22, 6).Value = 2189
$ws.Cells.Item(23, 6).Value = 1890
$ws.Cells.Item(24, 6).Value = 482
$ws.Cells.Item(26, 6).Value = 277
$ws.Cells.Item(27, 6).Value = 159
$ws.Cells.Item(28, 6).Value = 212
$ws.Cells.Item(30, 6).Value = 325
$ws.Cells.Item(31, 6).Value = 63
$ws.Cells.Item(32, 6).Value = 296
$ws.Cells.Item(33, 6).Value = 490
$ws.Cells.Item(34, 6).Value = 24
$ws.Cells.Item(37, 6).Value = 11
$ws.Cells.Item(38, 6).Value = 60
$ws.Cells.Item(39, 6).Value = 578
$ws.Cells.Item(41, 6).Value = 35
$ws.Cells.Item(42, 6).Value = 1396
$ws.Cells.Item(44, 6).Value = 302
$ws.Cells.Item(45, 6).Value = 9
$ws.Cells.Item(46, 6).Value = 176
$ws.Cells.Item(47, 6).Value = 642
$ws.Cells.Item(48, 6).Value = 297
